$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D; existing D:K shift to E:L
$ws.Columns.Item(4).Insert(-4161, 0)

# Copy number formats from the (shifted) old column D -- now column E -- into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the latest fiscal-period figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 17100800
$ws.Cells.Item(9, 4).Value = 7694800
$ws.Cells.Item(10, 4).Value = 9406100
$ws.Cells.Item(12, 4).Value = 129400
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 347000
$ws.Cells.Item(15, 4).Value = 3435400
$ws.Cells.Item(17, 4).Value = 12631500
$ws.Cells.Item(18, 4).Value = 4469300
$ws.Cells.Item(20, 4).Value = 16700
$ws.Cells.Item(21, 4).Value = 7921500
$ws.Cells.Item(22, 4).Value = 245100
$ws.Cells.Item(23, 4).Value = 4241000
$ws.Cells.Item(24, 4).Value = 874000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 3367000
$ws.Cells.Item(27, 4).Value = 3367000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 52000
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -16700
$ws.Cells.Item(33, 4).Value = 3419000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 3419000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 1555600
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 2343100
$ws.Cells.Item(44, 4).Value = 859400
$ws.Cells.Item(45, 4).Value = 299300
$ws.Cells.Item(46, 4).Value = 5057400
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 28075500
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 801600
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 33934500
$ws.Cells.Item(57, 4).Value = 2239800
$ws.Cells.Item(58, 4).Value = 913100
$ws.Cells.Item(59, 4).Value = 575400
$ws.Cells.Item(60, 4).Value = 3728400
$ws.Cells.Item(61, 4).Value = 5170200
$ws.Cells.Item(62, 4).Value = 5671800
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 14570300
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 13543100
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 19364200
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 3419000
$ws.Cells.Item(83, 4).Value = 3435400
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 7768600
$ws.Cells.Item(91, 4).Value = -6076500
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -6170200
$ws.Cells.Item(96, 4).Value = -438000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -839100
$ws.Cells.Item(101, 4).Value = -37900
$ws.Cells.Item(102, 4).Value = 721400
